$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet to allow edits, re-protect at the end to restore state.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer (shared string used by A81).
[void]$ws.Cells.Replace("2021-07-13", "2021-07-14")

# Update Weight (D) and Percent Change (E) for rows 2-78 (holdings data refresh).
$data = New-Object 'object[,]' 77,2
$data[0,0] = 0.08560757273342377
$data[0,1] = 0.02410052183466105
$data[1,0] = 0.04988219587819016
$data[1,1] = 0.001174755803076133
$data[2,0] = 0.04319589592597512
$data[2,1] = 0.005445227418321652
$data[3,0] = 0.03616139074821971
$data[3,1] = 0.01501959975242406
$data[4,0] = 0.03454692304491565
$data[4,1] = 0.007032271490440944
$data[5,0] = 0.02984920455215596
$data[5,1] = 0.006793879600637709
$data[6,0] = 0.02744744307049728
$data[6,1] = -0.003405075489881115
$data[7,0] = 0.02777421610841187
$data[7,1] = -0.003223380251423658
$data[8,0] = 0.02663066873026078
$data[8,1] = 0.007373877491421377
$data[9,0] = 0.0268510042794178
$data[9,1] = 0.0006870491240125176
$data[10,0] = 0.0219879776053082
$data[10,1] = -0.02508780732563975
$data[11,0] = 0.02325075848297233
$data[11,1] = -0.001252382248843054
$data[12,0] = 0.02025151825176289
$data[12,1] = 0.01444111027756945
$data[13,0] = 0.02081698436663211
$data[13,1] = -0.009079180006689946
$data[14,0] = 0.02145877439017711
$data[14,1] = -0.007155680295329869
$data[15,0] = 0.01879593811663815
$data[15,1] = -0.002193910934913945
$data[16,0] = 0.01779793970888239
$data[16,1] = 0.006899985773225259
$data[17,0] = 0.01734967239216933
$data[17,1] = 0.01847749602564441
$data[18,0] = 0.01751195112927677
$data[18,1] = -0.01266721576869545
$data[19,0] = 0.01486471423280721
$data[19,1] = -0.02217477003942181
$data[20,0] = 0.01433560144884753
$data[20,1] = 0.006844346317615546
$data[21,0] = 0.01303809500135935
$data[21,1] = -0.0003537318712415471
$data[22,0] = 0.01309859345503719
$data[22,1] = -0.01273766621101047
$data[23,0] = 0.0116583866190533
$data[23,1] = -0.01706484641638206
$data[24,0] = 0.0110063778731382
$data[24,1] = 0.004683263495193435
$data[25,0] = 0.01176464324546803
$data[25,1] = -0.0008647526807332717
$data[26,0] = 0.01044859840783941
$data[26,1] = 0.02253725917848048
$data[27,0] = 0.009932055556706667
$data[27,1] = 0.000509878903760308
$data[28,0] = 0.01028518928108925
$data[28,1] = 0.006638238009407882
$data[29,0] = 0.01014456880953614
$data[29,1] = 0.004100552683187653
$data[30,0] = 0.01028564143694634
$data[30,1] = -0.0061543871988744
$data[31,0] = 0.01062968682860154
$data[31,1] = 0.01389686459170791
$data[32,0] = 0.0104475132337824
$data[32,1] = -0.004154764996104898
$data[33,0] = 0.01003966865069266
$data[33,1] = 0.007566204287515754
$data[34,0] = 0.01095957974193001
$data[34,1] = -0.0003094250881860994
$data[35,0] = 0.008153997864309632
$data[35,1] = -0.005057226510513657
$data[36,0] = 0.009474790338437639
$data[36,1] = 0.0001288493750808417
$data[37,0] = 0.008473445977339461
$data[37,1] = -0.01086973922230927
$data[38,0] = 0.009656692639742526
$data[38,1] = -0.01247366203118427
$data[39,0] = 0.008458660480812818
$data[39,1] = -0.01885884730106813
$data[40,0] = 0.008991028786943493
$data[40,1] = 0.01032949790794979
$data[41,0] = 0.00906332850849122
$data[41,1] = -0.004664574675600108
$data[42,0] = 0.008123070403685087
$data[42,1] = 0.009306882194464849
$data[43,0] = 0.008936950946436248
$data[43,1] = 0.01820371157387735
$data[44,0] = 0.007870405710746561
$data[44,1] = -0.003171247357294016
$data[45,0] = 0.009301750291931601
$data[45,1] = -0.006844254326268673
$data[46,0] = 0.008211783382844964
$data[46,1] = -0.03055931811424228
$data[47,0] = 0.009068528300847686
$data[47,1] = -0.02267627965417174
$data[48,0] = 0.006810959322013091
$data[48,1] = -0.01649704911938277
$data[49,0] = 0.007642157434090568
$data[49,1] = 0.001230652719269143
$data[50,0] = 0.008108284907158441
$data[50,1] = 0.002091175240485077
$data[51,0] = 0.008511381853748808
$data[51,1] = 0.004462388440289011
$data[52,0] = 0.006580936333618044
$data[52,1] = -0.001447998309809906
$data[53,0] = 0.00656543869161649
$data[53,1] = -0.00499989669634926
$data[54,0] = 0.005564501270789688
$data[54,1] = -0.002925259616791132
$data[55,0] = 0.006061782282410647
$data[55,1] = 0.005012531328320913
$data[56,0] = 0.005606280471984247
$data[56,1] = -0.01572707476409374
$data[57,0] = 0.005867400479450245
$data[57,1] = 0.00300543289793076
$data[58,0] = 0.005133777601331491
$data[58,1] = -0.001497269684692659
$data[59,0] = 0.005074680830810614
$data[59,1] = -0.005978633735175709
$data[60,0] = 0.005254412784001495
$data[60,1] = 0.007710312543026232
$data[61,0] = 0.004878942560278959
$data[61,1] = 0.006524317912218303
$data[62,0] = 0.004859409427252932
$data[62,1] = 0.006699419383653371
$data[63,0] = 0.004496599567528748
$data[63,1] = -0.01053817070227658
$data[64,0] = 0.004102229228980102
$data[64,1] = -0.004563190265193917
$data[65,0] = 0.003894147103550055
$data[65,1] = -0.008360039013515452
$data[66,0] = 0.004056923212100287
$data[66,1] = -0.01076635014043059
$data[67,0] = 0.004162637251486521
$data[67,1] = -0.0004344898003518782
$data[68,0] = 0.004085544677853703
$data[68,1] = -0.001460871874896186
$data[69,0] = 0.0032945884370567
$data[69,1] = 0
$data[70,0] = 0.003310052167368972
$data[70,1] = -0.006515859355790687
$data[71,0] = 0.003662462442380224
$data[71,1] = -0.02017283950617288
$data[72,0] = 0.002735588150943777
$data[72,1] = 0.006313945224045936
$data[73,0] = 0.002444716288081654
$data[73,1] = 0.01344603092402186
$data[74,0] = 0.001896341664610205
$data[74,1] = -0.03633762517882677
$data[75,0] = 0.001446446586811647
$data[75,1] = -0.01181619256017519
$data[76,0] = 0.9999999999999999
$data[76,1] = 0.001800857651490517

$ws.Range("D2:E78").Value = $data

$ws.Protect()
